# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right before the existing "2022-Q1" sheet.
# 2) Copy header / index-column formatting from "2022-Q1" so the new sheet matches
#    the look of the other quarterly sheets.
# 3) Fill in the three funds reported for 2022-Q4.
# 4) Update the "总计" (totals) sheet with a new summary row for 2022-Q4, pushing
#    the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$q1 = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Add($q1)
$q4.Name = "2022-Q4"

# ---- formatting: reuse the header / row-index styles already used elsewhere ----
$q1.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$q4.Range("A2:A4").PasteSpecial(-4122)

$q1.Range("A1").Copy()
$q4.Range("A1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- header row ----
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# ---- data rows: B,D,E,F,G must stay TEXT (leading zeros / fixed decimals matter) ----
$q4.Range("B2:G4").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "004332"
$q4.Range("C2").Value = "恒生前海沪港深新兴产业精选混合"
$q4.Range("D2").Value = "0.49"
$q4.Range("E2").Value = "75.81"
$q4.Range("F2").Value = "3.53"
$q4.Range("G2").Value = "0.0173"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "013383"
$q4.Range("C3").Value = "恒生前海高端制造混合A"
$q4.Range("D3").Value = "0.11"
$q4.Range("E3").Value = "84.98"
$q4.Range("F3").Value = "5.67"
$q4.Range("G3").Value = "0.0062"
$q4.Range("H3").Value = 5

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "013384"
$q4.Range("C4").Value = "恒生前海高端制造混合C"
$q4.Range("D4").Value = "0.04"
$q4.Range("E4").Value = "84.98"
$q4.Range("F4").Value = "5.67"
$q4.Range("G4").Value = "0.0023"
$q4.Range("H4").Value = 5

# ---- update the "总计" summary sheet: insert a 2022-Q4 row, push the rest down ----
$total = $wb.Worksheets.Item("总计")
$total.Range("A2").EntireRow.Insert(-4121)
$total.Range("A2:D2").ClearFormats()

# the row-index cell (column A) carries the bordered/bold style - reuse it
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.03

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.29

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.44

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q3"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0
